# Update handback/handoff timestamps produced when the handback report was
# regenerated for 02f302fe-a608-412f-be9c-dea8f05e5926.md (zh-cn / de-de).

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the md file.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-30 12:57:38"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-30 12:57:33"
$wsZhCn.Range("K3").Value = "2016-08-30 12:58:48"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-30 12:57:38"
$wsDeDe.Range("K3").Value = "2016-08-30 12:58:55"
